$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format on price cells whose new values would
# otherwise be auto-coerced to numeric/scientific by COM Value assignment,
# so they keep rendering exactly like the scraped text (e.g. trailing zeros).
$textCells = @("D4", "D5", "D7", "D8", "D9", "D10", "D11", "D12", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D45", "D46", "D48", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the refreshed coin data scraped by the GitHub Actions job.
$ws.Range("D2").Value = "28.176.71"
$ws.Range("E2").Value = "  -3.30%  "
$ws.Range("D3").Value = "1.917.34"
$ws.Range("E3").Value = "  -3.98%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -1.20%  "
$ws.Range("D5").Value = "327.74"
$ws.Range("E5").Value = "  -0.97%  "
$ws.Range("E6").Value = "  -1.09%  "
$ws.Range("D7").Value = "0.4681"
$ws.Range("E7").Value = "  -5.88%  "
$ws.Range("D8").Value = "0.4019"
$ws.Range("E8").Value = "  -4.14%  "
$ws.Range("D9").Value = "52.91"
$ws.Range("E9").Value = "  -3.51%  "
$ws.Range("D10").Value = "0.08414"
$ws.Range("E10").Value = "  -5.15%  "
$ws.Range("D11").Value = "1.047"
$ws.Range("E11").Value = "  -4.36%  "
$ws.Range("D12").Value = "22.19"
$ws.Range("E12").Value = "  -3.42%  "
$ws.Range("D13").Value = "1.926.87"
$ws.Range("E13").Value = "  -4.62%  "
$ws.Range("E14").Value = "  -6.77%  "
$ws.Range("D15").Value = "6.080"
$ws.Range("E15").Value = "  -5.42%  "
$ws.Range("D16").Value = "1.003"
$ws.Range("E16").Value = "  -1.14%  "
$ws.Range("D17").Value = "89.69"
$ws.Range("E17").Value = "  -3.11%  "
$ws.Range("D18").Value = "0.00001065"
$ws.Range("E18").Value = "  -3.75%  "
$ws.Range("D19").Value = "0.06600"
$ws.Range("E19").Value = "  -2.28%  "
$ws.Range("D20").Value = "17.96"
$ws.Range("E20").Value = "  -7.88%  "
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  -1.15%  "
$ws.Range("D22").Value = "5.748"
$ws.Range("E22").Value = "  -3.77%  "
$ws.Range("D23").Value = "28.173.06"
$ws.Range("E23").Value = "  -3.40%  "
$ws.Range("E24").Value = "  -6.55%  "
$ws.Range("E25").Value = "  +0.43%  "
$ws.Range("D26").Value = "2.129.66"
$ws.Range("E26").Value = "  -5.56%  "
$ws.Range("D27").Value = "153.37"
$ws.Range("E27").Value = "  -2.37%  "
$ws.Range("D28").Value = "20.05"
$ws.Range("E28").Value = "  -3.68%  "
$ws.Range("D29").Value = "5.778"
$ws.Range("E29").Value = "  -8.01%  "
$ws.Range("D30").Value = "2.136"
$ws.Range("E30").Value = "  -5.04%  "
$ws.Range("D31").Value = "123.59"
$ws.Range("E31").Value = "  -2.82%  "
$ws.Range("D32").Value = "0.9806"
$ws.Range("E32").Value = "  -6.03%  "
$ws.Range("D33").Value = "0.09675"
$ws.Range("E33").Value = "  -1.95%  "
$ws.Range("D34").Value = "1.442"
$ws.Range("E34").Value = "  -5.87%  "
$ws.Range("D35").Value = "3.644"
$ws.Range("E35").Value = "  -2.30%  "
$ws.Range("D36").Value = "5.551"
$ws.Range("E36").Value = "  -4.73%  "
$ws.Range("D37").Value = "8.856"
$ws.Range("E37").Value = "  -2.45%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.02301"
$ws.Range("E38").Value = "  -4.80%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "1.262"
$ws.Range("E39").Value = "  -3.86%  "
$ws.Range("D40").Value = "0.06189"
$ws.Range("E40").Value = "  -3.03%  "
$ws.Range("D41").Value = "0.6175"
$ws.Range("E41").Value = "  -4.70%  "
$ws.Range("E42").Value = "  -4.40%  "
$ws.Range("E43").Value = "  -1.13%  "
$ws.Range("E44").Value = "  -3.83%  "
$ws.Range("D45").Value = "1.309"
$ws.Range("E45").Value = "  -3.41%  "
$ws.Range("D46").Value = "0.5861"
$ws.Range("E46").Value = "  -5.31%  "
$ws.Range("E47").Value = "  -3.62%  "
$ws.Range("D48").Value = "2.028"
$ws.Range("E48").Value = "  -6.76%  "
$ws.Range("D49").Value = "3.437"
$ws.Range("E49").Value = "  -1.89%  "
$ws.Range("E50").Value = "  -0.67%  "
$ws.Range("D51").Value = "111.37"
$ws.Range("E51").Value = "  -1.52%  "
